$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A190").Value = 189
$ws.Range("B190").Value = 1
$ws.Range("C190").Value = "2024-06-18 23:13:38"
$ws.Range("D190").Value = 200
$ws.Range("E190").Value = 13

$ws.Range("A191").Value = 190
$ws.Range("B191").Value = 2
$ws.Range("C191").Value = "2024-06-18 23:13:38"
$ws.Range("D191").Value = 200
$ws.Range("E191").Value = 2
